$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A23").Value = "2025-04-28 21:49:50"
$ws.Range("B23").Value = 53
